$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new "2022-Q3" row at the top of the
#    data block (row 2), pushing every existing row down by one.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Shift existing data rows (B:D) down by one, working from the bottom up
# so we never overwrite a row before it has been copied.
for ($r = 8; $r -ge 2; $r--) {
    $src = $total.Range("B" + $r + ":D" + $r)
    $dst = $total.Range("B" + ($r + 1) + ":D" + ($r + 1))
    $src.Copy($dst)
}

# Column A is just the zero-based row index; extend it to the new last row
# by copying an existing A-column cell (keeps its style) then fixing the
# value.
$total.Range("A8").Copy($total.Range("A9"))
$total.Range("A9").Value = 7

# Write the new "2022-Q3" row.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01

# ---------------------------------------------------------------------
# 2) Add a new worksheet "2022-Q3" right after "总计" holding the detail
#    rows for the quarter (same layout as the other quarterly sheets).
#    Duplicating an existing quarterly sheet (rather than Worksheets.Add)
#    keeps all of its formatting (sheetPr, header/row styles, margins).
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($null, $total)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# "2022-Q2" had two detail rows; "2022-Q3" only has one, so drop row 3.
$newSheet.Rows.Item(3).Delete()

# Overwrite row 2 with the "2022-Q3" detail values. Columns B-G are text
# in the source data (note the leading zero in the fund code and the
# fixed decimal formatting), so force a text number format before
# assigning, then drop back to the default style so no stray formatting
# is left on the cells.
$textCells = @("B2", "C2", "D2", "E2", "F2", "G2")
foreach ($addr in $textCells) {
    $newSheet.Range($addr).NumberFormat = "@"
}

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "008135"
$newSheet.Range("C2").Value = "华宸未来价值先锋混合"
$newSheet.Range("D2").Value = "0.28"
$newSheet.Range("E2").Value = "87.82"
$newSheet.Range("F2").Value = "3.72"
$newSheet.Range("G2").Value = "0.0104"

foreach ($addr in $textCells) {
    $newSheet.Range($addr).Style = "Normal"
}

$newSheet.Range("H2").Value = 9
